$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '91.961.81'
$ws.Cells.Item(2, 5).Value = '  +1.55%  '

$ws.Cells.Item(3, 4).Value = '3.099.45'
$ws.Cells.Item(3, 5).Value = '  +0.29%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = "'" + '239.43'
$ws.Cells.Item(5, 5).Value = '  -0.63%  '

$ws.Cells.Item(6, 4).Value = "'" + '613.78'
$ws.Cells.Item(6, 5).Value = '  -0.92%  '

$ws.Cells.Item(7, 5).Value = '  -5.34%  '

$ws.Cells.Item(8, 4).Value = "'" + '0.390'
$ws.Cells.Item(8, 5).Value = '  +7.19%  '

$ws.Cells.Item(9, 5).Value = '  -0.04%  '

$ws.Cells.Item(10, 4).Value = '3.099.13'
$ws.Cells.Item(10, 5).Value = '  +0.40%  '

$ws.Cells.Item(11, 4).Value = "'" + '0.727'
$ws.Cells.Item(11, 5).Value = '  -1.72%  '

$ws.Cells.Item(12, 5).Value = '  -1.58%  '

$ws.Cells.Item(13, 4).Value = "'" + '0.0000250'
$ws.Cells.Item(13, 5).Value = '  +1.18%  '

$ws.Cells.Item(14, 4).Value = '91.881.21'
$ws.Cells.Item(14, 5).Value = '  +1.57%  '

$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).Value = "'" + '34.31'
$ws.Cells.Item(15, 5).Value = '  -1.74%  '

$ws.Cells.Item(16, 2).Value = 'Toncoin'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(16, 4).Value = "'" + '5.50'
$ws.Cells.Item(16, 5).Value = '  +0.20%  '

$ws.Cells.Item(17, 4).Value = '3.685.21'
$ws.Cells.Item(17, 5).Value = '  +0.64%  '

$ws.Cells.Item(18, 4).Value = '3.088.78'
$ws.Cells.Item(18, 5).Value = '  -0.10%  '

$ws.Cells.Item(19, 4).Value = "'" + '3.63'
$ws.Cells.Item(19, 5).Value = '  -2.27%  '

$ws.Cells.Item(20, 4).Value = "'" + '14.71'
$ws.Cells.Item(20, 5).Value = '  +0.09%  '

$ws.Cells.Item(21, 4).Value = "'" + '5.80'
$ws.Cells.Item(21, 5).Value = '  +0.08%  '

$ws.Cells.Item(22, 4).Value = "'" + '445.86'
$ws.Cells.Item(22, 5).Value = '  +1.56%  '

$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(23, 4).Value = "'" + '9.29'
$ws.Cells.Item(23, 5).Value = '  +1.42%  '

$ws.Cells.Item(24, 2).Value = 'PEPE'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(24, 4).Value = "'" + '0.0000201'
$ws.Cells.Item(24, 5).Value = '  -3.27%  '

$ws.Cells.Item(25, 4).Value = "'" + '5.60'
$ws.Cells.Item(25, 5).Value = '  -4.91%  '

$ws.Cells.Item(26, 4).Value = "'" + '86.88'
$ws.Cells.Item(26, 5).Value = '  -2.57%  '

$ws.Cells.Item(27, 4).Value = "'" + '11.64'
$ws.Cells.Item(27, 5).Value = '  -1.59%  '

$ws.Cells.Item(28, 4).Value = '3.275.78'
$ws.Cells.Item(28, 5).Value = '  +0.49%  '

$ws.Cells.Item(29, 5).Value = '  -0.10%  '

$ws.Cells.Item(30, 5).Value = '  +12.58%  '

$ws.Cells.Item(31, 4).Value = "'" + '0.232'
$ws.Cells.Item(31, 5).Value = '  -6.06%  '

$ws.Cells.Item(32, 5).Value = '  -4.52%  '

$ws.Cells.Item(33, 4).Value = "'" + '9.19'
$ws.Cells.Item(33, 5).Value = '  -0.18%  '

$ws.Cells.Item(34, 5).Value = '  +57.14%  '

$ws.Cells.Item(35, 4).Value = "'" + '0.165'
$ws.Cells.Item(35, 5).Value = '  -3.39%  '

$ws.Cells.Item(36, 4).Value = "'" + '7.91'
$ws.Cells.Item(36, 5).Value = '  +0.08%  '

$ws.Cells.Item(37, 2).Value = 'MantraDAO'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(37, 4).Value = "'" + '4.23'
$ws.Cells.Item(37, 5).Value = '  -3.66%  '

$ws.Cells.Item(38, 2).Value = 'EthereumClassic'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(38, 4).Value = "'" + '26.08'
$ws.Cells.Item(38, 5).Value = '  -0.70%  '

$ws.Cells.Item(39, 4).Value = "'" + '1.93'
$ws.Cells.Item(39, 5).Value = '  +2.03%  '

$ws.Cells.Item(40, 4).Value = "'" + '480.02'
$ws.Cells.Item(40, 5).Value = '  -2.68%  '

$ws.Cells.Item(41, 5).Value = '  +0.12%  '

$ws.Cells.Item(42, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(42, 4).Value = "'" + '0.430'
$ws.Cells.Item(42, 5).Value = '  +2.24%  '

$ws.Cells.Item(43, 2).Value = 'dogwifhat'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(43, 4).Value = "'" + '3.41'
$ws.Cells.Item(43, 5).Value = '  -2.20%  '

$ws.Cells.Item(44, 4).Value = "'" + '22.18'
$ws.Cells.Item(44, 5).Value = '  +0.12%  '

$ws.Cells.Item(46, 4).Value = "'" + '158.84'
$ws.Cells.Item(46, 5).Value = '  +2.70%  '

$ws.Cells.Item(47, 4).Value = "'" + '1.89'
$ws.Cells.Item(47, 5).Value = '  -1.06%  '

$ws.Cells.Item(48, 4).Value = "'" + '0.693'
$ws.Cells.Item(48, 5).Value = '  +0.03%  '

$ws.Cells.Item(49, 4).Value = "'" + '1.36'
$ws.Cells.Item(49, 5).Value = '  +0.44%  '

$ws.Cells.Item(50, 4).Value = "'" + '0.0329'
$ws.Cells.Item(50, 5).Value = '  +3.60%  '

$ws.Cells.Item(51, 4).Value = "'" + '44.02'
$ws.Cells.Item(51, 5).Value = '  -0.49%  '
